# tracking sheets by city
# Update the confirmation summary so that Sofala's cities are broken out
# (adding Buzi and Muchungue), Maputo Provincia no longer lists "Ponta De
# Ouro", and the confirmado counts for several cities are refreshed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# We currently have 38 data rows (rows 2-38, plus header row 1). The new
# data needs 39 data rows (rows 2-39), so insert one extra row at the
# bottom of the table before rewriting the values.
$ws.Rows.Item(39).Insert()

# Rebuild rows 18 through 39 (Maputo Cidade .. Zambézia) with the updated
# province/city pairs and confirmado counts.
$rows = @(
    @(18, "Maputo Cidade",     "Maputo",      817),
    @(19, "Maputo Provincia",  "Boane",        19),
    @(20, "Maputo Provincia",  "Kathembe",      1),
    @(21, "Maputo Provincia",  "Manhiça",      34),
    @(22, "Maputo Provincia",  "Marracuene",   13),
    @(23, "Maputo Provincia",  "Matola",      235),
    @(24, "Maputo Provincia",  "Matutuine",    26),
    @(25, "Maputo Provincia",  "Namaacha",     25),
    @(26, "Nampula",           "Liuto",         1),
    @(27, "Nampula",           "Mama",          1),
    @(28, "Nampula",           "Nampula",       4),
    @(29, "Sofala",            "Beira",       164),
    @(30, "Sofala",            "Buzi",          1),
    @(31, "Sofala",            "Caia",          1),
    @(32, "Sofala",            "Dondo",         2),
    @(33, "Sofala",            "Muchungue",     1),
    @(34, "Sofala",            "Nhamantada",    3),
    @(35, "Tete",              "Moatize",       1),
    @(36, "Tete",              "Tete",          3),
    @(37, "Zambézia",          "Nicoadala",    13),
    @(38, "Zambézia",          "Quelimane",    26),
    @(39, "Zambézia",          "Zavala",        1)
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
}
